$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: new build entry ---
# Copy formatting (date number format) from the row above so we reuse the
# existing style instead of creating a new custom number format.
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A23").Value = 45260.868055555555
$ws.Range("B23").Value = 1774592
$ws.Range("C23").Value = 612864
$ws.Range("D23").Value = 392192
$ws.Range("F23").Value = 144108396
$ws.Range("I23").Value = "Remove all q2 enemies except infantry, add Zombie, remove all legacy OGL and more soft code, fix Draw_Fill usage"

# --- Column I width (widened, no longer "best fit") ---
$ws.Columns.Item(9).ColumnWidth = 102.5

# --- Selection moved to F23 ---
$ws.Range("F23").Select() | Out-Null
